$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: header row ----
$ws.Range("A2").Value = "Unique ID"
$ws.Range("B2").Value = "Course Code"
$ws.Range("C2").Value = "Title"
$ws.Range("D2").Value = "Category"
$ws.Range("E2").Value = "Start Date"
$ws.Range("F2").Value = "Start Time"
$ws.Range("G2").Value = "End Date"
$ws.Range("H2").Value = "End Time"
$ws.Range("I2").Value = "Timezone"
$ws.Range("J2").Value = "Location"
$ws.Range("K2").Value = "Description"
$ws.Range("L2").Value = "Link"
$ws.Range("M2").Value = "TRANSPARENT"

# ---- Course code column ----
$ws.Range("B3").Value = "BBUS23"
$ws.Range("B4").Value = "MBAA22"

# ---- Title / Category columns (entry order matches original authoring) ----
$ws.Range("C3").Value = "Welcome to course"
$ws.Range("D4").Value = "Orientation"
$ws.Range("C4").Value = "Trimester Orientation"
$ws.Range("D3").Value = "Lecture"

# ---- Start/End dates ----
$ws.Range("E3").Value = 45884
$ws.Range("G3").Value = 45884
$ws.Range("E4").Value = 45885
$ws.Range("G4").Value = 45885

# ---- Start/End times ----
$ws.Range("F4").Value = "9:30am"
$ws.Range("F3").Value = "8:00am"
$ws.Range("H3").Value = "11:30am"
$ws.Range("H4").Value = "2:00pm"

# ---- Timezone / Location ----
$ws.Range("I3").Value = "Sydney, Australia"
$ws.Range("J3").Value = "Sydney, Australia"
$ws.Range("I4").Value = "Sydney, Australia"
$ws.Range("J4").Value = "Sydney, Australia"

# ---- Description ----
$ws.Range("K3").Value = "First lecture"
$ws.Range("K4").Value = "Frist Orientation"

# ---- Link ----
$ws.Range("L3").Value = "https:www.torrens.edu.au"
$ws.Range("L4").Value = "https:www.torrens.edu.au"

# ---- Transparent flag ----
$ws.Range("M3").Value = "Transparent"
$ws.Range("M4").Value = "Transparent"

# ---- Unique ID formulas (concat of course code + start date + location) ----
$ws.Range("A3").Formula = "=B3&E3&J3"
$ws.Range("A4").Formula = "=B4&E4&J4"

# ---- Number formats: register date format (style 1) before time format (style 2) ----
$ws.Range("E3").NumberFormat = "d-mmm"
$ws.Range("G3").NumberFormat = "d-mmm"
$ws.Range("E4").NumberFormat = "d-mmm"
$ws.Range("G4").NumberFormat = "d-mmm"
$ws.Range("I4").NumberFormat = "d-mmm"
$ws.Range("K4").NumberFormat = "d-mmm"
$ws.Range("M4").NumberFormat = "d-mmm"

$ws.Range("F3").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

# ---- Final selection (matches last-saved cursor position) ----
[void]$ws.Range("M5").Select()
